$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for new columns J (10) and K (11) ---
$ws.Columns.Item(10).ColumnWidth = 391.5
$ws.Columns.Item(11).ColumnWidth = 440.8

# --- Copy header formatting (bold, border, wrap) from A1 onto J1:K1 ---
$ws.Range("A1").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)

# --- Copy data-cell formatting (wrap text) from G2 onto J2:K15 ---
$ws.Range("G2").Copy()
$ws.Range("J2:K15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header values ---
$ws.Range("J1").Value = 'Onkelos'
$ws.Range("K1").Value = 'Jonathan'

# --- Translation text values (Onkelos / Jonathan columns) ---
$ws.Range("J3").Value = 'A flawless lamb, a yearling male must be in your possession. You may take it from sheep or goats.'
$ws.Range("K3").Value = 'The lamb shall be perfect, a male, the son of a year he shall be to you; from the sheep or from the young goats ye may take.'
$ws.Range("J4").Value = 'He said, Please [<b>Now</b>] take your son, your only one, who you love—Yitzchok—and go to the land of Moriah [<b>worship</b>]. Sacrifice him [<b>before me</b>] as a burnt-offering on one of the mountains which I will designate to you.'
$ws.Range("K4").Value = 'And He said, Take now thy son, thy only one whom thou lovest, Izhak, and go into the land of worship, and offer him there, a whole burnt offering, upon one of the mountains that I will tell thee.'
$ws.Range("J9").Value = 'Anything that has a blemish on it you shall not bring, for it will not be accepted favorably for you.'
$ws.Range("K9").Value = 'But anything that hath a blemish you shall not offer; for that will not be acceptable from you.'
$ws.Range("J10").Value = 'I see it but not now, I perceive it but not in the near future; a star [<b>king</b>] has gone forth from Yaakov, and a staff has arisen [<b>the Moshiach will be magnified by</b>] from Yisroel, which will smash the corners [<b>kill the leaders</b>] of Moav, and impale all of the sons of Sheis [<b>will rule over mankind</b>].'
$ws.Range("K10").Value = 'I shall see Him, but not now; I shall behold Him, but it is not near. When the mighty King of Jakob''s house shall reign, and the Meshiha, the Power-sceptre of Israel, be anointed, He will slay the princes of the Moabaee, and bring to nothing all the children of Sheth, the armies of Gog who will do battle against Israel and all their carcases shall fall before Him.'
$ws.Range("J11").Value = 'But if you will not do this, behold, you will have sinned against [<b>before</b>] Adonoy and you must realize that your sin will find you!'
$ws.Range("K11").Value = 'But if you will not perform this, behold, ye will have sinned before the Lord your God, and know that your sin will meet you.'
